$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in I1 from "DoB_Excel" to "BirthDate_Excel"
$ws.Range("I1").Value = "BirthDate_Excel"

# Update the selection to J19 (as recorded in the sheetView)
$ws.Range("J19").Select()
